$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row of mail-log data ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A6").Value = "Kun jij dit even regelen?"
$logs.Range("B6").Value = "mailmind.test@zohomail.eu"
$logs.Range("C6").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D6").Value = "Overig"
$logs.Range("E6").Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Range("F6").Value = "2025-08-01 23:36:06"
$logs.Range("G6").Value = "Ja"
$logs.Range("H6").Value = "Ja"
$logs.Range("I6").Value = "Nee"
$logs.Range("J6").Value = "Nee"

# Extend the conditional-formatting ranges so they cover the new row 6 too
$logs.Range("D2:D5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D6"))
$logs.Range("G2:G5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G6"))
$logs.Range("H2:H5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H6"))
$logs.Range("I2:I5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I6"))
$logs.Range("J2:J5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J6"))

# --- Sheet "Dashboard": bump the "Overig" count from 4 to 5 ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 5
